$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 28, shifting the existing rows 28-62 down to 29-63.
$ws.Rows.Item(28).EntireRow.Insert()

# Populate the newly inserted row 28 with a new weekly price entry
# (same product dimensions as the rest of the block, new date/volume/prices).
$ws.Range("A28").Value = 4
$ws.Range("B28").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C28").Value = "Los Lagos"
$ws.Range("D28").Value = 44915
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100103
$ws.Range("H28").Value = "Frutos de hueso (carozo)"
$ws.Range("I28").Value = 100103003
$ws.Range("J28").Value = "Damasco"
$ws.Range("K28").Value = "Castle Brite"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 800
$ws.Range("N28").Value = 19000
$ws.Range("O28").Value = 20000
$ws.Range("P28").Value = 19500
$ws.Range("Q28").Value = "$/caja 16 kilos"
$ws.Range("R28").Value = "Región de O'Higgins"
$ws.Range("S28").Value = 1219
$ws.Range("T28").Value = 16
